$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 82
$ws.Range("B82").Value = 5574442
$ws.Range("F82").Value = "FK Qarabag"
$ws.Range("G82").Value = "FK Sumqayit"
$ws.Range("H82").Value = 1
$ws.Range("I82").Value = 2
$ws.Range("J82").Value = "A"
$ws.Range("K82").Value = 1.125
$ws.Range("L82").Value = 7.5
$ws.Range("M82").Value = 15
$ws.Range("N82").Value = 1.2
$ws.Range("O82").Value = 6
$ws.Range("P82").Value = 11
$ws.Range("Q82").Value = -2.25
$ws.Range("R82").Value = 1.975
$ws.Range("S82").Value = 1.825
$ws.Range("T82").Value = 3.5
$ws.Range("U82").Value = 1.825
$ws.Range("V82").Value = 1.975
$ws.Range("X82").Value = -1
$ws.Range("Y82").Value = 10
$ws.Range("AA82").Value = 0.825
$ws.Range("AC82").Value = 0.9750000000000001

# Row 83
$ws.Range("B83").Value = 5573343
$ws.Range("F83").Value = "Shamakhi FK"
$ws.Range("G83").Value = "FK Gabala"
$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = "D"
$ws.Range("K83").Value = 3.5
$ws.Range("L83").Value = 3.1
$ws.Range("M83").Value = 2
$ws.Range("N83").Value = 3.3
$ws.Range("O83").Value = 3.2
$ws.Range("P83").Value = 2.05
$ws.Range("Q83").Value = 0.25
$ws.Range("R83").Value = 2
$ws.Range("S83").Value = 1.8
$ws.Range("T83").Value = 2.5
$ws.Range("U83").Value = 1.975
$ws.Range("V83").Value = 1.825
$ws.Range("X83").Value = 2.2
$ws.Range("Y83").Value = -1
$ws.Range("Z83").Value = 0.5
$ws.Range("AA83").Value = -0.5
$ws.Range("AC83").Value = 0.825

# Row 85
$ws.Range("B85").Value = 5579144
$ws.Range("F85").Value = "Sabah"
$ws.Range("G85").Value = "Zira IK"
$ws.Range("K85").Value = 1.45
$ws.Range("L85").Value = 4.2
$ws.Range("M85").Value = 5.5
$ws.Range("N85").Value = 1.5
$ws.Range("O85").Value = 4
$ws.Range("P85").Value = 5.25
$ws.Range("Q85").Value = -1
$ws.Range("R85").Value = 1.85
$ws.Range("S85").Value = 1.95
$ws.Range("U85").Value = 1.8
$ws.Range("V85").Value = 2
$ws.Range("X85").Value = 3
$ws.Range("Z85").Value = -1
$ws.Range("AA85").Value = 0.95
$ws.Range("AC85").Value = 1

